$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.491899999999998
$ws.Range("B3").Value = 5.85469999999999
$ws.Range("E4").Value = 13.11530000000001
$ws.Range("B5").Value = 5.419499999999998
$ws.Range("E6").Value = 12.29129999999999
$ws.Range("D7").Value = -6.643499999999997
$ws.Range("A9").Value = -20.15309999999998
$ws.Range("D9").Value = -8.644400000000008
$ws.Range("E10").Value = 11.9476
$ws.Range("B11").Value = 4.895800000000005
$ws.Range("E11").Value = 13.1258
$ws.Range("B12").Value = 5.220000000000002
$ws.Range("A13").Value = -21.80490000000002
$ws.Range("A16").Value = -19.9065
$ws.Range("A18").Value = -22.80770000000002
$ws.Range("A20").Value = -22.08240000000004
$ws.Range("B21").Value = 4.828500000000002
$ws.Range("D21").Value = -7.6825
$ws.Range("E21").Value = 13.12
$ws.Range("E25").Value = 13.174
